$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Shape 1: Title placeholder ("Assembly scaffolding" -> empty)
$s.Shapes.Item(1).TextFrame.TextRange.Text = ""

# Shape 2: Content placeholder ("Connecting contigs..." -> empty)
$s.Shapes.Item(2).TextFrame.TextRange.Text = ""
